$wb = $excel.ActiveWorkbook

# --- "daily" sheet: keep only the first 4 rows (Tickers, GBP_USD, USD_JPY, EUR_USD) ---
$daily = $wb.Worksheets.Item("daily")
$daily.Rows("5:8").Delete() | Out-Null
$daily.Range("D6").Select() | Out-Null

# --- "minute" sheet: append GBP_USD as a third row ---
$minute = $wb.Worksheets.Item("minute")
$minute.Range("A3").Value = "GBP_USD"
$minute.Range("A2:A3").Select() | Out-Null

# --- Add a new sheet ("Sheet1") after "minute" with the remaining tickers ---
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $minute)
$newSheet.Name = "Sheet1"

$values = @("USD_JPY", "EUR_USD", "AUD_USD", "EUR_JPY", "GBP_JPY", "NZD_USD", "EUR_USD", "GBP_USD")
for ($i = 0; $i -lt $values.Length; $i++) {
    $newSheet.Cells.Item($i + 1, 1).Value = $values[$i]
}

$newSheet.Range("M21").Select() | Out-Null
$newSheet.Activate()
